$d = $word.ActiveDocument

$p3 = $d.Paragraphs(3)
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$rng = $d.Range($p3.Range.Start, $pLast.Range.End)

$frag = @'
<w:p><w:r><w:t>Crea</w:t></w:r><w:r><w:t>ț</w:t></w:r><w:r><w:t>i o noua metoda statică</w:t></w:r><w:r><w:t xml:space="preserve"> ce utilizează </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>String</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> și metode precum</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>indexOf</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>lastIndexOf</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>length</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>isEmpty</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>split</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Exemple de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>string</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>Ala bala portocala</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>To</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>be</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, or </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>not</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>to</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>be</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>that</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>is</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>the</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>question</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Whether</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 'tis </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nobler</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>the</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mind</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>to</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>suffer</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>slings</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>and</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>arrows</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>outrageous</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fortune</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Or </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>to</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>take</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>arms</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>against</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sea</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>troubles</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:r><w:t>TASK 2</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>În proiectul de săptămâna trecută adăugați la model două clase noi</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Invoice</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> și </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>InvoiceItem</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>Adăugați câmpuri precum</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>invoiceDate</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>array</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> of  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>invoiceItems</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>InvoiceItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>has</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">a </w:t></w:r><w:r><w:t>Product</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">ca </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>field</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">și TVA as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>double</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, tot </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>field</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>TASK 3</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Faceți clasa Product abstractă și adăugații două, trei metode abstracte (e.g. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hiddenInformation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>computePrice</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>Creați un număr de clase ce extind Product e.g.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>FoodProduct</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ClotheProduct</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>CleaningProduct</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>FurnitureProduct</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Toate aceste clase faceți-le imutabile, excepție </w:t></w:r><w:r><w:t xml:space="preserve">fiind </w:t></w:r><w:r><w:t>Product.</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t xml:space="preserve">TASK 4 </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Pentru toate clasele din model suprascrieți</w:t></w:r><w:r><w:t xml:space="preserve"> metoda</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>toString</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> folosind </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>StringBuilder</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Creați-vă un obiect (instanță) din fiecare clasă și folosiți-o într-un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>System.out.println</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pentru a verifica rezultatul.</w:t></w:r></w:p><w:p/><w:p/>
'@

$xml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>$frag</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

[void]$rng.InsertXML($xml)

Write-Output "Paragraph count after edit: $($d.Paragraphs.Count)"
